$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.440.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.644.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.68%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.23"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.49"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.69%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.642.91"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.68%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.59"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.66%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.121.64"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.346.02"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.31%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.639.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.93"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.35"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.71"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.20%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.09"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.65"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.97%  "

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "SuiNetwork"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.53"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "541.21"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +17.97%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.37"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.76"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.81"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +13.61%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0803"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "175.26"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.73%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.84"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +8.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.402"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.00"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.79"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +9.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.30"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.24%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.72"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.28"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0563"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.631"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0959"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0238"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.70"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.83%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.35"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.58%  "
